$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B column values (row -> value)
$bValues = @{
    2  = 0.92439
    3  = 0.92793
    4  = 0.92416
    5  = 0.92472
    6  = 0.91569
    7  = 0.91899
    8  = 0.85976
    9  = 0.92779
    10 = 0.86849
    11 = 0.93372
}

# New C column text values (row -> text), matches new shared-string table
$cValues = @{
    2  = "21:45:00"
    3  = "19:25:00"
    4  = "20:50:00"
    5  = "19:06:00"
    6  = "19:13:00"
    7  = "18:59:00"
    8  = "20:36:00"
    9  = "21:39:00"
    10 = "21:45:00"
    11 = "19:01:00"
}

foreach ($r in 2..11) {
    $ws.Range("B$r").Value = $bValues[$r]
    $ws.Range("C$r").Value = $cValues[$r]
}
